$d = $word.ActiveDocument

# Fix the typo "nekoretni" -> "nekorektni" in the sentence
# "... namerno uvredljivi, nekoretni ili cak greskom ostavljeni utisci."
# Only this specific occurrence (followed by " ili") should change; the
# other occurrences of "nekoretnih"/"nekoretni" elsewhere in the document
# (title, ToC, summary) must stay untouched.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "nekoretni ili"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "nekorektni ili"
$find.Forward = $true
$find.Wrap = 0
$find.Format = $false
$find.MatchCase = $true
$find.MatchWholeWord = $false
$find.MatchWildcards = $false
$find.MatchSoundsLike = $false
$find.MatchAllWordForms = $false
$find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$find.Replacement.Text, [ref]2) | Out-Null
